# Regenerate merged AHB files
#
# The workbook has one sheet ("AHB-Diff") with a header row (row 1) and
# data rows 2..61 across columns A..U. Columns A..J (plus K="diff") hold
# the "_old" comparison columns, and columns L..U hold the "_new" ones.
# This edit:
#   1. Renames the "_old" headers to "_FV2410" and the "_new" headers to
#      "_FV2504" (column K, "diff", is left untouched).
#   2. Freezes the header row (row 1) in the sheet view.
#   3. Turns the used range A1:U61 into a native Excel Table ("Table1").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AHB-Diff")

# --- 1. Rename header row ------------------------------------------------

$renames = @{
    "A1" = "Segmentname_FV2410"
    "B1" = "Segmentgruppe_FV2410"
    "C1" = "Segment_FV2410"
    "D1" = "Datenelement_FV2410"
    "E1" = "Segment ID_FV2410"
    "F1" = "Code_FV2410"
    "G1" = "Qualifier_FV2410"
    "H1" = "Beschreibung_FV2410"
    "I1" = "Bedingungsausdruck_FV2410"
    "J1" = "Bedingung_FV2410"
    "L1" = "Segmentname_FV2504"
    "M1" = "Segmentgruppe_FV2504"
    "N1" = "Segment_FV2504"
    "O1" = "Datenelement_FV2504"
    "P1" = "Segment ID_FV2504"
    "Q1" = "Code_FV2504"
    "R1" = "Qualifier_FV2504"
    "S1" = "Beschreibung_FV2504"
    "T1" = "Bedingungsausdruck_FV2504"
    "U1" = "Bedingung_FV2504"
}

foreach ($addr in $renames.Keys) {
    $ws.Range($addr).Value = $renames[$addr]
}

# --- 2. Freeze the header row ---------------------------------------------

$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true

# --- 3. Convert A1:U61 into a native table ---------------------------------

$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $ws.Range("A1:U61"), $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"

Write-Output "Renamed headers, froze header row, and created Table1 over A1:U61"
